# Update import product detail table format.xlsx (import maks 10000 row berikutnya)
# The first data row (row 3) is removed; subsequent rows shift up by one,
# which naturally renumbers the Size Code / Product code pairing sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 entirely, shifting rows 4:65 up to 3:64.
$ws.Rows("3:3").Delete()

# Leave the selection where the editor ended up after performing the edit.
$ws.Range("C14").Select()
